# Update (Analyze PO & Forecast)
$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$newForecast = @(47, 48, 48, 48, 47, 47, 48, 48, 47, 39, 30, 29, 33, 38, 38, 36)
for ($i = 0; $i -lt $newForecast.Length; $i++) {
    $row = $i + 2
    $wsForecast.Range("D$row").Value = $newForecast[$i]
}

# --- Sheet: Summary ---
# These cells hold numeric-looking / date-looking values that must stay
# stored as literal TEXT (matching the source data's inlineStr cells), so
# a leading apostrophe is used to force text entry, just as typing an
# apostrophe-prefixed value into Excel does.
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").Value = "'671"
$wsSummary.Range("B10").Value = "'381"
$wsSummary.Range("B11").Value = "'191"
$wsSummary.Range("B12").Value = "'48"
$wsSummary.Range("B13").Value = "'2025-02-02"
$wsSummary.Range("B14").Value = "'29"
$wsSummary.Range("B15").Value = "'2025-04-13"
